# Added Validation for the order Period in the Product Dependency Feature
#
# 1. Rename the test-data "Web Data 49" company (and every string derived
#    from it) to "Web Data 1" everywhere it appears in the workbook.
# 2. Leave the ConfigOrderPeriods sheet's "order period" validation work
#    as the last thing the user touched: land on the GeneratePayInvoice
#    sheet with cell N28 selected (was C1).

$wb = $excel.ActiveWorkbook

# --- 1. Global rename: "Web Data 49" -> "Web Data 1" -------------------
# A plain substring replace on every sheet handles the base string plus
# every string that embeds it ("Web Data 49 Child", "Web Data 49
# Reseller", the two "Successfully created ..." messages, and "Working
# as admin Web Data 49 Child X") in one pass.
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("Web Data 49", "Web Data 1")
}

# --- 2. Update the active sheet / selection -----------------------------
$target = $wb.Worksheets.Item("GeneratePayInvoice")
$target.Activate()
$target.Range("N28").Select()
